$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.861.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.838.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.17%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'231.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.67%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.90%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'39.59"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -5.57%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.327"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.24%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0684"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.40%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0985"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.103.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.75%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'WrappedEther"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'1.840.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'Polygon"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.673"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.73%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'34.848.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'69.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.00%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'240.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.33%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.53%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.63%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.38%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'171.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.58%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.16%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +2.42%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'17.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.83%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.98%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.09%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0552"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.60%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.99%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +4.12%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.88%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +11.56%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.86%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'91.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.70%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +5.28%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.340.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.35%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.64%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'14.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.09%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.74%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.10%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'6.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.30%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0522"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.67%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.017.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.00%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0681"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.06%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D51").Value = "'3.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +14.47%  "
$ws.Range("E51").Style = "Normal"
